# Add a new worksheet "EmpData" at the end of the workbook to hold the
# employee-creation test data (per commit: "added employee creation testcase").

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "EmpData"

# Header row
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "Id"
$ws.Range("D1").Value = "Username"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "ConfirmPassword"

# Data row
$ws.Range("A2").Value = "Test1"
$ws.Range("B2").Value = "Test2"
$ws.Range("C2").Value = "Test3"
$ws.Range("D2").Value = "Test1"
$ws.Range("E2").Value = "Test1"
$ws.Range("F2").Value = "Test1"

# Size the first few columns to fit their content, like the other sheets
# (A/B best-fit to the header text, C widened to comfortably fit "ConfirmPassword").
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 8.333333333333334
$ws.Columns.Item(3).ColumnWidth = 18.333333333333334

# Match the selection left behind on the new sheet.
$ws.Range("I13").Select() | Out-Null
